# actualizacion de produtos en pagina
#
# Adds a "Resultado" column and swaps the single demo product row for two
# generic product rows (producto1 / producto2), formatting the header row
# with a bold, centered, bordered style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: Producto | Precio | Resultado -----------------------------
$ws.Range("A1").Value = "Producto"
$ws.Range("B1").Value = "Precio"
$ws.Range("C1").Value = "Resultado"

# --- Data rows ---------------------------------------------------------
$ws.Range("A2").Value = "producto1"
$ws.Range("B2").Value = 1150000

$ws.Range("A3").Value = "producto2"
$ws.Range("B3").Value = 600000

# --- Header styling ------------------------------------------------------
# Build the combined font/border/alignment format once on a scratch cell
# (well outside the used range) and copy it onto the header range in a
# single paste so only one new font / border / cell style gets registered,
# instead of one per individual property assignment.
$scratch = $ws.Range("Z1")
$scratch.Font.Bold = $true
$scratch.Borders.LineStyle = 1
$scratch.HorizontalAlignment = -4108
$scratch.VerticalAlignment = -4160
$scratch.Copy()
$ws.Range("A1:C1").PasteSpecial(-4122)
$scratch.Clear()

# --- Selection -------------------------------------------------------------
[void]$ws.Range("D5").Select()
